# Cluster Data mining.xlsx - relabel clusters 1-5 -> 0-4, move the footnote
# down one row, clear the stale "General" number-format override from the
# Amount/Sum rows, and update the view (zoom + selection).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Relabel the five "Cluster 1..5" header blocks to "Cluster 0..4" ---
$blocks = @("B2","H2","N2","T2","Z2")
foreach ($start in $blocks) {
    $rng = $ws.Range($start).Resize(1, 5)
    $rng.Cells.Item(1,1).Value = "Cluster 0"
    $rng.Cells.Item(1,2).Value = "Cluster 1"
    $rng.Cells.Item(1,3).Value = "Cluster 2"
    $rng.Cells.Item(1,4).Value = "Cluster 3"
    $rng.Cells.Item(1,5).Value = "Cluster 4"
}

# --- 2. Move the footnote from row 8 down to row 9 ---
$note = $ws.Range("A8").Value2
$ws.Range("A8").ClearContents()
$ws.Range("A9").Value = $note

# --- 3. Clear the leftover explicit "General" number format on the Amount
#        row (B5:F5) and the Sum row totals, matching the surrounding
#        default-styled cells (copy the default format from A5) ---
$ws.Range("A5").Copy()
$ws.Range("B5:F5").PasteSpecial(-4122)
$ws.Range("B6").PasteSpecial(-4122)
$ws.Range("H6").PasteSpecial(-4122)
$ws.Range("N6").PasteSpecial(-4122)
$ws.Range("T6").PasteSpecial(-4122)
$ws.Range("Z6").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- 4. Update the view: zoom to 145% and move the selection to H14 ---
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.Zoom = 145
$ws.Range("H14").Select()
